$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.627.82'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.13%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.843.22'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.74%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.08'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.67%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4246'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -2.34%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3634'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.42%  '
$ws.Range('B9').NumberFormat = '@'
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').NumberFormat = '@'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07273'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.78%  '
$ws.Range('B10').NumberFormat = '@'
$ws.Range('B10').Value = 'Polygon'
$ws.Range('C10').NumberFormat = '@'
$ws.Range('C10').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8875'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -5.41%  '
$ws.Range('B11').NumberFormat = '@'
$ws.Range('B11').Value = 'Solana'
$ws.Range('C11').NumberFormat = '@'
$ws.Range('C11').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.63'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -3.05%  '
$ws.Range('B12').NumberFormat = '@'
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').NumberFormat = '@'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.815.98'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.83%  '
$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.560'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.95%  '
$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').NumberFormat = '@'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.337'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.65%  '
$ws.Range('B15').NumberFormat = '@'
$ws.Range('B15').Value = 'TRON'
$ws.Range('C15').NumberFormat = '@'
$ws.Range('C15').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06890'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.50%  '
$ws.Range('B16').NumberFormat = '@'
$ws.Range('B16').Value = 'BinanceUSD'
$ws.Range('C16').NumberFormat = '@'
$ws.Range('C16').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.003'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.11%  '
$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').NumberFormat = '@'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '78.89'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -3.10%  '
$ws.Range('B18').NumberFormat = '@'
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').NumberFormat = '@'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008871'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.93%  '
$ws.Range('B19').NumberFormat = '@'
$ws.Range('B19').Value = 'Dai'
$ws.Range('C19').NumberFormat = '@'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.001'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.04%  '
$ws.Range('B20').NumberFormat = '@'
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').NumberFormat = '@'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.43'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -2.73%  '
$ws.Range('B21').NumberFormat = '@'
$ws.Range('B21').Value = 'WrappedBTC'
$ws.Range('C21').NumberFormat = '@'
$ws.Range('C21').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '27.616.48'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.03%  '
$ws.Range('B22').NumberFormat = '@'
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').NumberFormat = '@'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.976'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -2.42%  '
$ws.Range('B23').NumberFormat = '@'
$ws.Range('B23').Value = 'Cosmos'
$ws.Range('C23').NumberFormat = '@'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.54'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -3.89%  '
$ws.Range('B24').NumberFormat = '@'
$ws.Range('B24').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C24').NumberFormat = '@'
$ws.Range('C24').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.036.86'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.97%  '
$ws.Range('B25').NumberFormat = '@'
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').NumberFormat = '@'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.949'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.91%  '
$ws.Range('B26').NumberFormat = '@'
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').NumberFormat = '@'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '155.43'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.07%  '
$ws.Range('B27').NumberFormat = '@'
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').NumberFormat = '@'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.98'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +3.42%  '
$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'BitcoinCash'
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '122.87'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +8.36%  '
$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.254'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -2.47%  '
$ws.Range('B30').NumberFormat = '@'
$ws.Range('B30').Value = 'LidoDAOToken'
$ws.Range('C30').NumberFormat = '@'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.843'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +6.71%  '
$ws.Range('B31').NumberFormat = '@'
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').NumberFormat = '@'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08940'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.27%  '
$ws.Range('B32').NumberFormat = '@'
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').NumberFormat = '@'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7741'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.61%  '
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.573'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -5.35%  '
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.924'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -2.28%  '
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.092'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -6.67%  '
$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value = 'Frax'
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.001'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.097'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -1.71%  '
$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05360'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.61%  '
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01933'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.74%  '
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.793'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -4.61%  '
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.863'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.81%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5072'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -3.07%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1650'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.74%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.255'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -5.86%  '
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.06586'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -2.08%  '
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.34'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.78%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4710'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -3.45%  '
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '104.68'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.35%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'PaxDollar'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.000'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.07%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.635'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.46%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.769'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -8.46%  '
